$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.82195472717285
$ws.Range("C3").Value = 17.77195930480957
$ws.Range("C4").Value = 17.63105392456055
$ws.Range("C5").Value = 17.46082305908203
$ws.Range("C6").Value = 18.08714866638184
